# Update the frequency table values in B2:X5 (rows correspond to A/C/G/T,
# columns correspond to position 1..23) with the refreshed run's statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 4,23
$arr[0,0] = 0.042332268370607
$arr[0,1] = 0.0375399361022364
$arr[0,2] = 0.952076677316294
$arr[0,3] = 0.0519169329073482
$arr[0,4] = 0.0167731629392971
$arr[0,5] = 0.0359424920127796
$arr[0,6] = 0.958466453674121
$arr[0,7] = 0.0311501597444089
$arr[0,8] = 0.969648562300319
$arr[0,9] = 0.0207667731629393
$arr[0,10] = 0.669329073482428
$arr[0,11] = 0.962460063897764
$arr[0,12] = 0.144568690095847
$arr[0,13] = 0.097444089456869
$arr[0,14] = 0.00479233226837061
$arr[0,15] = 0.982428115015974
$arr[0,16] = 0.00399361022364217
$arr[0,17] = 0.975239616613419
$arr[0,18] = 0.976038338658147
$arr[0,19] = 0.232428115015974
$arr[0,20] = 0.958466453674121
$arr[0,21] = 0.0535143769968051
$arr[0,22] = 0.012779552715655
$arr[1,0] = 0.0167731629392971
$arr[1,1] = 0.00319488817891374
$arr[1,2] = 0.00878594249201278
$arr[1,3] = 0.10223642172524
$arr[1,4] = 0.94888178913738
$arr[1,5] = 0.938498402555911
$arr[1,6] = 0.0295527156549521
$arr[1,7] = 0.00638977635782748
$arr[1,8] = 0.00559105431309904
$arr[1,9] = 0.0527156549520767
$arr[1,10] = 0.154952076677316
$arr[1,11] = 0.00878594249201278
$arr[1,12] = 0
$arr[1,13] = 0.00319488817891374
$arr[1,14] = 0.00159744408945687
$arr[1,15] = 0
$arr[1,16] = 0.0207667731629393
$arr[1,17] = 0.0175718849840256
$arr[1,18] = 0
$arr[1,19] = 0.00479233226837061
$arr[1,20] = 0.0103833865814696
$arr[1,21] = 0.00878594249201278
$arr[1,22] = 0.0103833865814696
$arr[2,0] = 0.133386581469649
$arr[2,1] = 0.935303514376997
$arr[2,2] = 0.0319488817891374
$arr[2,3] = 0.662140575079872
$arr[2,4] = 0.0135782747603834
$arr[2,5] = 0.00718849840255591
$arr[2,6] = 0.00638977635782748
$arr[2,7] = 0.957667731629393
$arr[2,8] = 0.0159744408945687
$arr[2,9] = 0.619009584664537
$arr[2,10] = 0.166932907348243
$arr[2,11] = 0.0023961661341853
$arr[2,12] = 0.833067092651757
$arr[2,13] = 0.89776357827476
$arr[2,14] = 0.00638977635782748
$arr[2,15] = 0.0175718849840256
$arr[2,16] = 0.975239616613419
$arr[2,17] = 0.000798722044728434
$arr[2,18] = 0.0215654952076677
$arr[2,19] = 0.756389776357827
$arr[2,20] = 0.00878594249201278
$arr[2,21] = 0.915335463258786
$arr[2,22] = 0.972044728434505
$arr[3,0] = 0.807507987220447
$arr[3,1] = 0.023961661341853
$arr[3,2] = 0.00559105431309904
$arr[3,3] = 0.18370607028754
$arr[3,4] = 0.0207667731629393
$arr[3,5] = 0.0175718849840256
$arr[3,6] = 0.00559105431309904
$arr[3,7] = 0
$arr[3,8] = 0.00878594249201278
$arr[3,9] = 0.307507987220447
$arr[3,10] = 0.00479233226837061
$arr[3,11] = 0.0223642172523962
$arr[3,12] = 0.0223642172523962
$arr[3,13] = 0.00159744408945687
$arr[3,14] = 0.987220447284345
$arr[3,15] = 0
$arr[3,16] = 0
$arr[3,17] = 0.00638977635782748
$arr[3,18] = 0.0023961661341853
$arr[3,19] = 0.00638977635782748
$arr[3,20] = 0.0223642172523962
$arr[3,21] = 0.0215654952076677
$arr[3,22] = 0.00479233226837061

$ws.Range("B2:X5").Value = $arr
